# Applies the "Updated cryptos list" price/volume refresh (GitHub Actions scrape).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '34.395.23'
$ws.Range('E2').Value = '  -0.45%  '

# Row 3
$ws.Range('D3').Value = '1.801.09'
$ws.Range('E3').Value = '  +0.19%  '

# Row 4
$ws.Range('E4').Value = '  +0.55%  '

# Row 5
$ws.Range('D5').Value = '''227.52'
$ws.Range('E5').Value = '  +0.31%  '

# Row 6
$ws.Range('D6').Value = '''0.578'
$ws.Range('E6').Value = '  +3.46%  '

# Row 7
$ws.Range('E7').Value = '  +0.62%  '

# Row 8
$ws.Range('D8').Value = '''34.84'
$ws.Range('E8').Value = '  +5.36%  '

# Row 9
$ws.Range('D9').Value = '''0.299'
$ws.Range('E9').Value = '  +0.26%  '

# Row 10
$ws.Range('D10').Value = '''0.0691'
$ws.Range('E10').Value = '  -0.66%  '

# Row 11
$ws.Range('D11').Value = '''0.0951'
$ws.Range('E11').Value = '  +0.09%  '

# Row 12
$ws.Range('D12').Value = '2.061.38'
$ws.Range('E12').Value = '  +0.45%  '

# Row 13
$ws.Range('D13').Value = '''11.18'
$ws.Range('E13').Value = '  +0.15%  '

# Row 14
$ws.Range('D14').Value = '1.812.96'
$ws.Range('E14').Value = '  +1.20%  '

# Row 15
$ws.Range('D15').Value = '''0.642'
$ws.Range('E15').Value = '  +0.33%  '

# Row 16
$ws.Range('D16').Value = '34.349.97'
$ws.Range('E16').Value = '  -0.21%  '

# Row 17
$ws.Range('D17').Value = '''4.33'
$ws.Range('E17').Value = '  +0.99%  '

# Row 18
$ws.Range('D18').Value = '''69.11'
$ws.Range('E18').Value = '  +0.19%  '

# Row 19
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').Value = '0.0₃0795'
$ws.Range('E19').Value = '  -0.84%  '

# Row 20
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').Value = '''244.63'
$ws.Range('E20').Value = '  -1.74%  '

# Row 21
$ws.Range('D21').Value = '''11.49'
$ws.Range('E21').Value = '  +0.92%  '

# Row 22
$ws.Range('E22').Value = '  +0.34%  '

# Row 23
$ws.Range('D23').Value = '''4.15'
$ws.Range('E23').Value = '  -0.81%  '

# Row 24
$ws.Range('D24').Value = '''170.97'
$ws.Range('E24').Value = '  +3.65%  '

# Row 25
$ws.Range('E25').Value = '  +2.14%  '

# Row 26
$ws.Range('D26').Value = '''7.51'
$ws.Range('E26').Value = '  +3.23%  '

# Row 27
$ws.Range('D27').Value = '''16.73'
$ws.Range('E27').Value = '  +0.89%  '

# Row 28
$ws.Range('E28').Value = '  +1.51%  '

# Row 29
$ws.Range('E29').Value = '  +0.29%  '

# Row 32
$ws.Range('E32').Value = '  +0.56%  '

# Row 33
$ws.Range('E33').Value = '  -0.42%  '

# Row 34
$ws.Range('D34').Value = '''1.83'
$ws.Range('E34').Value = '  -0.18%  '

# Row 35
$ws.Range('D35').Value = '1.399.62'
$ws.Range('E35').Value = '  -1.65%  '

# Row 36
$ws.Range('E36').Value = '  -1.62%  '

# Row 37
$ws.Range('D37').Value = '''0.676'
$ws.Range('E37').Value = '  +0.27%  '

# Row 38
$ws.Range('E38').Value = '  +0.00%  '

# Row 39
$ws.Range('E39').Value = '  -1.73%  '

# Row 40
$ws.Range('D40').Value = '''82.81'
$ws.Range('E40').Value = '  -2.78%  '

# Row 41
$ws.Range('E41').Value = '  +3.54%  '

# Row 42
$ws.Range('D42').Value = '''0.947'
$ws.Range('E42').Value = '  +0.78%  '

# Row 43
$ws.Range('D43').Value = '''2.40'
$ws.Range('E43').Value = '  +0.57%  '

# Row 44
$ws.Range('D44').Value = '''13.67'
$ws.Range('E44').Value = '  +1.13%  '

# Row 45
$ws.Range('E45').Value = '  +2.79%  '

# Row 46
$ws.Range('D46').Value = '''0.0510'
$ws.Range('E46').Value = '  -2.04%  '

# Row 47
$ws.Range('E47').Value = '  -1.49%  '

# Row 48
$ws.Range('D48').Value = '1.962.62'
$ws.Range('E48').Value = '  +0.59%  '

# Row 49
$ws.Range('D49').Value = '''104.57'
$ws.Range('E49').Value = '  -1.21%  '

# Row 50
$ws.Range('E50').Value = '  +0.41%  '

# Row 51
$ws.Range('D51').Value = '0.0₆0128'
$ws.Range('E51').Value = '  +0.14%  '
